# Insert a new daily price record as row 48 in the "Uva" sheet.
# This shifts the existing rows 48-126 down to 49-127 (dimension A1:T126 -> A1:T127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 48 (pushes rows 48..126 down to 49..127)
$ws.Rows.Item(48).Insert()

# Populate the new row 48 with the new record's data
$ws.Range("A48").Value = 8
$ws.Range("B48").Value = "Terminal La Palmera de La Serena"
$ws.Range("C48").Value = "Coquimbo"
$ws.Range("D48").Value = (Get-Date -Year 2022 -Month 6 -Day 9 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("E48").Value = 4
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100109
$ws.Range("H48").Value = "Uva"
$ws.Range("I48").Value = 100109001
$ws.Range("J48").Value = "Uva"
$ws.Range("K48").Value = "Red Globe"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 300
$ws.Range("N48").Value = 9500
$ws.Range("O48").Value = 10000
$ws.Range("P48").Value = 9750
$ws.Range("Q48").Value = "`$/bandeja 18 kilos"
$ws.Range("R48").Value = "Provincia de Limarí"
$ws.Range("S48").Value = 542
$ws.Range("T48").Value = 18
